$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.373.69"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").Value = "1.684.92"
$ws.Range("E3").Value = "  +3.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.31%  "
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0625"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0902"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("D12").Value = "1.926.91"
$ws.Range("E12").Value = "  +3.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +16.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.619"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.66%  "
$ws.Range("D15").Value = "1.688.17"
$ws.Range("E15").Value = "  +3.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.53%  "
$ws.Range("D17").Value = "30.360.34"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "0.0₃0721"
$ws.Range("E20").Value = "  +2.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.72%  "
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0502"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.54%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.73%  "
$ws.Range("D34").Value = "1.510.03"
$ws.Range("E34").Value = "  +5.73%  "
$ws.Range("E35").Value = "  +5.46%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  +5.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.586"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "79.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.43%  "
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.853"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.76%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0505"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "51.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.35%  "
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.48%  "
$ws.Range("D51").Value = "0.0₆0114"
$ws.Range("E51").Value = "  +6.45%  "
